$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Cell C10 on sheet "Rules" changes its value from 18 to 100
$ws.Range("C10").Value = 100
